# TC_Template/Main.rvl.xlsx edit:
# 1. Use of Global Object - the four "Functions" rows that used to call
#    CrmLaunchSales / CrmChangeArea / CrmOpenEntity / CrmClickButton (Object=Functions)
#    are refactored to call the shared "Crm" global object with short action names.
# 2. Update Range maps (removed optional params) - the optional fromRow/fromCol/toRow/toCol
#    Param rows under the "Map Range" block are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# --- 2. Remove the four optional Range-map Param rows (old rows 9-12: fromRow/fromCol/toRow/toCol) ---
$ws.Rows("9:12").Delete()

# --- 1. Use of Global Object: rename Object/Action pairs for the Crm function calls ---
# After the row deletion above, the four "Functions" rows are now at rows 11-14.
$ws.Range("C11").Value = "Crm"
$ws.Range("D11").Value = "LaunchSales"

$ws.Range("C12").Value = "Crm"
$ws.Range("D12").Value = "ChangeArea"

$ws.Range("C13").Value = "Crm"
$ws.Range("D13").Value = "OpenEntity"

$ws.Range("C14").Value = "Crm"
$ws.Range("D14").Value = "ClickButton"
